# Apply the "corrected code" results update to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (Hydrogen): update Iron & steel demand, clear Non-metallic minerals value
$ws.Range("B3").Value = 1987965.510719053
$ws.Range("D3").ClearContents()

# Row 4 (Methanol): update Chemicals value
$ws.Range("C4").Value = 11.25061418147256

# Row 5 (Ammonia): update Chemicals value
$ws.Range("C5").Value = 0

# Row 7: rename "Other" -> "Biogas" and update its value
$ws.Range("A7").Value = "Biogas"
$ws.Range("D7").Value = 180.0719470142306

# Row 8 (new): "Other" row, matching formatting of row 7, with new value
$ws.Range("A7:D7").Copy()
$ws.Range("A8:D8").PasteSpecial(-4122)
$ws.Range("A8").Value = "Other"
$ws.Range("B8").ClearContents()
$ws.Range("C8").ClearContents()
$ws.Range("D8").Value = 688.9084094431955

$excel.CutCopyMode = 0
